$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75: new departure (Hurghada, FlyEgypt)
$ws.Cells.Item(75, 1).Value = 74.0
$ws.Cells.Item(75, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(75, 3).Value = "8:30 PM"
$ws.Cells.Item(75, 4).Value = "FT1518"
$ws.Cells.Item(75, 5).Value = "Hurghada"
$ws.Cells.Item(75, 6).Value = "(HRG)"
$ws.Cells.Item(75, 7).Value = "FlyEgypt "
$ws.Cells.Item(75, 8).Value = "B738"
$ws.Cells.Item(75, 9).Value = "(SU-TMN)"
$ws.Cells.Item(75, 10).Value = "8:15 PM"
$ws.Cells.Item(75, 12).Value = "0 hours, -15 minutes"

# Row 76: new departure (Cologne, Ryanair)
$ws.Cells.Item(76, 1).Value = 75.0
$ws.Cells.Item(76, 2).Value = "Sunday, Jan 08"
$ws.Cells.Item(76, 3).Value = "10:05 PM"
$ws.Cells.Item(76, 4).Value = "FR6868"
$ws.Cells.Item(76, 5).Value = "Cologne"
$ws.Cells.Item(76, 6).Value = "(CGN)"
$ws.Cells.Item(76, 7).Value = "Ryanair "
$ws.Cells.Item(76, 8).Value = "B738"
$ws.Cells.Item(76, 9).Value = "(9H-QBA)"
$ws.Cells.Item(76, 10).Value = "10:20 PM"
$ws.Cells.Item(76, 12).Value = "0 hours, 15 minutes"
